# Update column C (date) from 45175 (2023-09-06) to 45177 (2023-09-08)
# for all data rows (2 through 74) on the "Avverkningsanmälningar" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$lastRow = 74
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45175) {
        $cell.Value2 = 45177
    }
}
